# Update the report from "through 12-18" to "through 12-19": rename the
# sheet, relabel the December row, and refresh the December / Total rows
# with one additional day's worth of "no arrest made" cases (arrest_made
# counts are unchanged; arrest_rate is recomputed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the only worksheet/tab to reflect the new "through" date.
$ws.Name = "Through 2021-12-19"

# Relabel the December row.
$ws.Range("A14").Value = "December (through 12-19)"

# --- Row 14 (December) updates ---------------------------------------
# 2015
$ws.Range("C14").Value = 21
$ws.Range("D14").Value = 0.125
# 2016
$ws.Range("F14").Value = 55
$ws.Range("G14").Value = 0.0833
# 2017
$ws.Range("I14").Value = 66
$ws.Range("J14").Value = 0.1081
# 2019
$ws.Range("O14").Value = 27
$ws.Range("P14").Value = 0.1
# 2020
$ws.Range("R14").Value = 85
$ws.Range("S14").Value = 0.0449
# 2021
$ws.Range("U14").Value = 131
$ws.Range("V14").Value = 0.015

# --- Row 15 (Total) updates -------------------------------------------
# 2015
$ws.Range("C15").Value = 279
$ws.Range("D15").Value = 0.1143
# 2016
$ws.Range("F15").Value = 559
$ws.Range("G15").Value = 0.1027
# 2017
$ws.Range("I15").Value = 824
$ws.Range("J15").Value = 0.0793
# 2019
$ws.Range("O15").Value = 507
$ws.Range("P15").Value = 0.1011
# 2020
$ws.Range("R15").Value = 1285
$ws.Range("S15").Value = 0.0503
# 2021
$ws.Range("U15").Value = 1673
$ws.Range("V15").Value = 0.058
